$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$t = $m.Theme
$tcs = $t.GetType().InvokeMember("ThemeColorScheme", [System.Reflection.BindingFlags]::GetProperty, $null, $t, $null)
$name = $tcs.GetType().InvokeMember("Name", [System.Reflection.BindingFlags]::GetProperty, $null, $tcs, $null)
Write-Host "SlideMaster Theme ColorScheme Name:" $name
